$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (names) for rows 8-15 shift down two slots (line7/line8 were
# inserted ahead of the "extr" block), and two brand-new rows (16,17) are
# appended at the bottom.
$ws.Range("B8").Value  = "line7"
$ws.Range("B9").Value  = "line8"
$ws.Range("B10").Value = "extr1"
$ws.Range("B11").Value = "extr2"
$ws.Range("B12").Value = "extr3"
$ws.Range("B13").Value = "extr4"
$ws.Range("B14").Value = "extr5"
$ws.Range("B15").Value = "extr6"

# Updated from_bus / to_bus / in_service values for the shifted rows.
$ws.Range("C8").Value  = 14
$ws.Range("D8").Value  = 11

$ws.Range("C9").Value  = 16
$ws.Range("E9").Value  = $true

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true

$ws.Range("D13").Value = 8

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# New rows 16 and 17 (contingencies extr7 / extr8), formatted like the rest
# of column A (bold, bordered, centered).
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
